# Update the NFL lookup table for the 2021-2022 season:
# Washington's old short code "WAS" is replaced with "WSH" and the
# placeholder team name "********" becomes "Football Team".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Team name first (so the new shared string for it is appended before WSH,
# matching the order the edits were made in).
$ws.Range("C33").Value = "Football Team"
$ws.Range("A33").Value = "WSH"

# Leave the selection where the author last left it.
[void]$ws.Range("F27").Select()
